# Updated IPS AIP hipo turnover
# Applies updated turnover figures across several location sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Bowling Green Ohio
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bowling Green Ohio")

# Professional Voluntary Turnover - ytd
$ws.Range("E2").Value = 0.0396
$ws.Range("E3").Value = 0.0396
$ws.Range("E4").Value = 0.0396

# Professional Voluntary Turnover - Commit/Forecast monthly/quarterly/FY
$ws.Range("M4").Value = 0.0204
$ws.Range("N4").Value = 0.04
$ws.Range("O4").Value = 0.0066
$ws.Range("P4").Value = 0.0066
$ws.Range("Q4").Value = 0.0066
$ws.Range("R4").Value = 0.0198
$ws.Range("S4").Value = 0.0066
$ws.Range("T4").Value = 0.0066
$ws.Range("U4").Value = 0.0066
$ws.Range("V4").Value = 0.0198
$ws.Range("W4").Value = 0.0792

# Internal Fill Rate - Commit/Forecast: clear Jun/Q2 values
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# Manufacturing Voluntary Turnover - ytd
$ws.Range("E8").Value = 0.0233
$ws.Range("E9").Value = 0.0233
$ws.Range("E10").Value = 0.0233

# Manufacturing Voluntary Turnover - Commit/Forecast monthly/quarterly/FY
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0.0229
$ws.Range("O10").Value = 0.00388333333333333
$ws.Range("P10").Value = 0.00388333333333333
$ws.Range("Q10").Value = 0.00388333333333333
$ws.Range("R10").Value = 0.01165
$ws.Range("S10").Value = 0.00388333333333333
$ws.Range("T10").Value = 0.00388333333333333
$ws.Range("U10").Value = 0.00388333333333333
$ws.Range("V10").Value = 0.01165
$ws.Range("W10").Value = 0.0466

# ---------------------------------------------------------------------------
# Sheet: Hyderabad India
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")

# Professional Voluntary Turnover - Commit/Forecast monthly/quarterly/FY
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.00926666666666667
$ws.Range("P4").Value = 0.00926666666666667
$ws.Range("Q4").Value = 0.00926666666666667
$ws.Range("R4").Value = 0.0278
$ws.Range("S4").Value = 0.00926666666666667
$ws.Range("T4").Value = 0.00926666666666667
$ws.Range("U4").Value = 0.00926666666666667
$ws.Range("V4").Value = 0.0278
$ws.Range("W4").Value = 0.1112

# ---------------------------------------------------------------------------
# Sheet: Langley Canada
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Langley Canada")

# Professional Voluntary Turnover - Commit/Forecast monthly/quarterly/FY
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0116
$ws.Range("O4").Value = 0.00193333333333333
$ws.Range("P4").Value = 0.00193333333333333
$ws.Range("Q4").Value = 0.00193333333333333
$ws.Range("R4").Value = 0.0058
$ws.Range("S4").Value = 0.00193333333333333
$ws.Range("T4").Value = 0.00193333333333333
$ws.Range("U4").Value = 0.00193333333333333
$ws.Range("V4").Value = 0.0058
$ws.Range("W4").Value = 0.0232

# Internal Fill Rate - Commit/Forecast: clear Jun, set Q2 to 1
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 1

# Manufacturing Voluntary Turnover - ytd
$ws.Range("E8").Value = 0.0412
$ws.Range("E9").Value = 0.0412
$ws.Range("E10").Value = 0.0412

$ws.Range("H10").Value = 0.0137
$ws.Range("J10").Value = 0.0139

$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0.0272
$ws.Range("O10").Value = 0.00686666666666667
$ws.Range("P10").Value = 0.00686666666666667
$ws.Range("Q10").Value = 0.00686666666666667
$ws.Range("R10").Value = 0.0206
$ws.Range("S10").Value = 0.00686666666666667
$ws.Range("T10").Value = 0.00686666666666667
$ws.Range("U10").Value = 0.00686666666666667
$ws.Range("V10").Value = 0.0206
$ws.Range("W10").Value = 0.0824

# ---------------------------------------------------------------------------
# Sheet: Las Vegas Nevada
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Las Vegas Nevada")

# Professional Voluntary Turnover - ytd
$ws.Range("E2").Value = 0.3738
$ws.Range("E3").Value = 0.3738
$ws.Range("E4").Value = 0.3738

# Professional Voluntary Turnover - Commit/Forecast monthly/quarterly/FY
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.1942
$ws.Range("O4").Value = 0.0623
$ws.Range("P4").Value = 0.0623
$ws.Range("Q4").Value = 0.0623
$ws.Range("R4").Value = 0.1869
$ws.Range("S4").Value = 0.0623
$ws.Range("T4").Value = 0.0623
$ws.Range("U4").Value = 0.0623
$ws.Range("V4").Value = 0.1869
$ws.Range("W4").Value = 0.7476

# Internal Fill Rate - Commit/Forecast: clear Jun value
$ws.Range("M7").ClearContents()

# ---------------------------------------------------------------------------
# Sheet: Apodaca Pmc Plant 2 Mexico
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Apodaca Pmc Plant 2 Mexico")

# Professional Voluntary Turnover - ytd / Commit-Forecast monthly/quarterly/FY
$ws.Range("E2").Value = 0.1818
$ws.Range("K2").Value = 0.2
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.2
$ws.Range("O2").Value = 0.0303
$ws.Range("P2").Value = 0.0303
$ws.Range("Q2").Value = 0.0303
$ws.Range("R2").Value = 0.0909
$ws.Range("S2").Value = 0.0303
$ws.Range("T2").Value = 0.0303
$ws.Range("U2").Value = 0.0303
$ws.Range("V2").Value = 0.0909
$ws.Range("W2").Value = 0.3636

# Manufacturing Voluntary Turnover - ytd / Commit-Forecast monthly/quarterly/FY
$ws.Range("E3").Value = 0.1667
$ws.Range("M3").Value = 0.0435
$ws.Range("N3").Value = 0.1688
$ws.Range("O3").Value = 0.0277833333333333
$ws.Range("P3").Value = 0.0277833333333333
$ws.Range("Q3").Value = 0.0277833333333333
$ws.Range("R3").Value = 0.08335
$ws.Range("S3").Value = 0.0277833333333333
$ws.Range("T3").Value = 0.0277833333333333
$ws.Range("U3").Value = 0.0277833333333333
$ws.Range("V3").Value = 0.08335
$ws.Range("W3").Value = 0.3334
